$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 260 currently only has A (date), D (3940) and E (30) populated.
# Fill in the missing B and C values to match the standard pattern used
# by all the other data rows.
$ws.Cells.Item(260, 2).Value = 187
$ws.Cells.Item(260, 3).Value = 628

# Append new daily rows for 17-09-2021 through 29-09-2021 (rows 261-273),
# each following the standard pattern (B=187, C=628, D=3940, E=30).
$startRow = 261
for ($day = 17; $day -le 29; $day++) {
    $dateText = "{0:D2}-09-2021" -f $day
    $row = $startRow + ($day - 17)

    $ws.Cells.Item($row, 1).Value = $dateText
    $ws.Cells.Item($row, 2).Value = 187
    $ws.Cells.Item($row, 3).Value = 628
    $ws.Cells.Item($row, 4).Value = 3940
    $ws.Cells.Item($row, 5).Value = 30
}

# Row 274 is the final new row (30-09-2021), matching the pattern of the
# last existing row: only A, D and E are populated (no B/C).
$lastRow = 274
$ws.Cells.Item($lastRow, 1).Value = "30-09-2021"
$ws.Cells.Item($lastRow, 4).Value = 3940
$ws.Cells.Item($lastRow, 5).Value = 30
